# Réception des produits Farnell
#
# 1. Drop the "last visible row" thick-bottom-border formatting that Excel
#    had put on row 53 (it will move to row 54 once the filter below makes
#    row 54 the new last visible row of the table).
# 2. Switch the "Supplier 1" (column F) AutoFilter on the BOM table from
#    "ETML" to "Farnell" - this shows/hides the relevant rows automatically
#    and moves the table's thick bottom border onto the new last visible
#    row.
# 3. Mark the "Reçu" (received) checkbox column (L) for the newly-received
#    Farnell line items.
# 4. Move the selection to the cell that is active once the update is done.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clear the explicit "last row" formatting on row 53 ---------------
$ws.Rows.Item(53).AutoFit()

# --- 2. Update the table's AutoFilter on "Supplier 1" (6th column of the
#        table == column F) from "ETML" to "Farnell" ---------------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Range.AutoFilter(6, @("Farnell"), 7)

# --- 3. Tick the "Reçu" column (L) for the received Farnell items --------
$ws.Range("L17").Value = 1
$ws.Range("L22").Value = 1
$ws.Range("L25").Value = 1
$ws.Range("L54").Value = 1

# --- 4. Move the selection to the cell that ends up active ---------------
$ws.Range("M57").Select()
